$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 6494.3335
$ws.Range("I113").Value = 2983
$ws.Range("K113").Value = 2983
$ws.Range("M113").Value = 271

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 927.7646999999999
$ws.Range("I2").Value = 812.6429000000001
$ws.Range("J2").Value = 1465
$ws.Range("K2").Value = 812.6429000000001
$ws.Range("L2").Value = 1465
$ws.Range("M2").Value = -699.6429000000001
$ws.Range("N2").Value = -1691

$ws.Range("H4").Value = 178.6923
$ws.Range("I4").Value = 185.18182
$ws.Range("K4").Value = 185.18182
$ws.Range("M4").Value = -69.18181999999999

$ws.Range("H5").Value = 60.27778
$ws.Range("J5").Value = 79.72727
$ws.Range("L5").Value = 79.72727
$ws.Range("N5").Value = -303.72727

$ws.Range("H32").Value = 3851221.5
$ws.Range("J32").Value = 7146384.5
$ws.Range("L32").Value = 7146384.5
$ws.Range("N32").Value = -7146958.5

$ws.Range("H44").Value = 12024.789
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 12024.789
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 12024.789
$ws.Range("M44").ClearContents()
$ws.Range("N44").Value = -13000.789

$ws.Range("H55").Value = 0
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 3365.8572
$ws.Range("I61").Value = 2283.6155
$ws.Range("K61").Value = 2283.6155
$ws.Range("M61").Value = -2071.6155

$ws.Range("H116").Value = 927.7646999999999
$ws.Range("I116").Value = 812.6429000000001
$ws.Range("J116").Value = 1465
$ws.Range("K116").Value = 812.6429000000001
$ws.Range("L116").Value = 1465
$ws.Range("M116").Value = 1481.3571
$ws.Range("N116").Value = -6053

$ws.Range("H132").Value = 3376.15
$ws.Range("I132").Value = 3346.889
$ws.Range("K132").Value = 10040.667
$ws.Range("M132").Value = -7510.667000000001

$ws.Range("H136").Value = 3365.8572
$ws.Range("I136").Value = 2283.6155
$ws.Range("K136").Value = 6850.8465
$ws.Range("M136").Value = -4300.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 927.7646999999999
$ws.Range("I3").Value = 812.6429000000001
$ws.Range("J3").Value = 1465
$ws.Range("K3").Value = 812.6429000000001
$ws.Range("L3").Value = 1465
$ws.Range("M3").Value = -698.6429000000001
$ws.Range("N3").Value = -1693

$ws.Range("H4").Value = 60.27778
$ws.Range("J4").Value = 79.72727
$ws.Range("L4").Value = 79.72727
$ws.Range("N4").Value = -309.72727

$ws.Range("H137").Value = 75000
$ws.Range("J137").Value = 130000
$ws.Range("L137").Value = 130000
$ws.Range("N137").Value = -140200

$ws.Range("H138").Value = 205000
$ws.Range("J138").Value = 205000
$ws.Range("L138").Value = 205000
$ws.Range("N138").Value = -215280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5526.6
$ws.Range("J31").Value = 7163.846
$ws.Range("L31").Value = 7163.846
$ws.Range("N31").Value = -7753.846

$ws.Range("H34").Value = 5526.6
$ws.Range("J34").Value = 7163.846
$ws.Range("L34").Value = 7163.846
$ws.Range("N34").Value = -7567.846

$ws.Range("H99").Value = 2400
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws.Range("H126").Value = 2400
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 3317.2666
$ws.Range("I17").Value = 27
$ws.Range("J17").Value = 7077.5713
$ws.Range("K17").Value = 81
$ws.Range("L17").Value = 21232.7139
$ws.Range("M17").Value = 88
$ws.Range("N17").Value = -21570.7139

$ws.Range("H103").Value = 659.6667
$ws.Range("I103").Value = 292.83334
$ws.Range("K103").Value = 878.5000200000001
$ws.Range("M103").Value = 0.4999799999999368

$ws.Range("H131").Value = 1700.1578
$ws.Range("I131").Value = 973.25
$ws.Range("K131").Value = 2919.75
$ws.Range("M131").Value = 2120.25

$ws.Range("H132").Value = 1787.5555
$ws.Range("J132").Value = 1784.875
$ws.Range("L132").Value = 16063.875
$ws.Range("N132").Value = -21123.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 26740.234
$ws.Range("J43").Value = 45099.9
$ws.Range("L43").Value = 45099.9
$ws.Range("N43").Value = -45401.9

$ws.Range("H46").Value = 56382.5
$ws.Range("J46").Value = 98765
$ws.Range("L46").Value = 98765
$ws.Range("N46").Value = -99077

$ws.Range("H57").Value = 55588.332
$ws.Range("J57").Value = 73382.5
$ws.Range("L57").Value = 73382.5
$ws.Range("N57").Value = -75022.5

$ws.Range("H102").Value = 1710.4546
$ws.Range("I102").Value = 1801.7778
$ws.Range("J102").Value = 1299.5
$ws.Range("K102").Value = 1801.7778
$ws.Range("L102").Value = 1299.5
$ws.Range("M102").Value = -179.7778000000001
$ws.Range("N102").Value = -4543.5

$ws.Range("H126").Value = 4333.3335
$ws.Range("J126").Value = 5000
$ws.Range("L126").Value = 15000
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2674.875
$ws.Range("I61").Value = 1398.091
$ws.Range("J61").Value = 5483.8
$ws.Range("K61").Value = 1398.091
$ws.Range("L61").Value = 5483.8
$ws.Range("M61").Value = -1196.091
$ws.Range("N61").Value = -5887.8

$ws.Range("H82").Value = 6800
$ws.Range("I82").Value = 800
$ws.Range("J82").Value = 8000
$ws.Range("K82").Value = 800
$ws.Range("L82").Value = 8000
$ws.Range("M82").Value = -439
$ws.Range("N82").Value = -8722

$ws.Range("H85").Value = 6800
$ws.Range("I85").Value = 800
$ws.Range("J85").Value = 8000
$ws.Range("K85").Value = 800
$ws.Range("L85").Value = 8000
$ws.Range("M85").Value = 448
$ws.Range("N85").Value = -10496

$ws.Range("H113").Value = 2674.875
$ws.Range("I113").Value = 1398.091
$ws.Range("J113").Value = 5483.8
$ws.Range("K113").Value = 1398.091
$ws.Range("L113").Value = 5483.8
$ws.Range("M113").Value = 771.9090000000001
$ws.Range("N113").Value = -9823.799999999999

$ws.Range("H122").Value = 3995.6667
$ws.Range("I122").Value = 3995.6667
$ws.Range("K122").Value = 11987.0001
$ws.Range("M122").Value = -9537.000100000001

$ws.Range("H131").Value = 157500
$ws.Range("J131").Value = 65000
$ws.Range("L131").Value = 65000
$ws.Range("N131").Value = -75080

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 10000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 10000
$ws.Range("M26").ClearContents()
$ws.Range("N26").Value = -10586

$ws.Range("H107").Value = 3374.5
$ws.Range("I107").Value = 4198.3335
$ws.Range("J107").Value = 903
$ws.Range("K107").Value = 12595.0005
$ws.Range("L107").Value = 2709
$ws.Range("M107").Value = -10675.0005
$ws.Range("N107").Value = -6549
